$d = $word.ActiveDocument

$oldUrl = "https://www.linkedin.com/in/jordan-alexis-799214175/"
$keepPrefix = "https://www.linkedin.com/in/jordan-alexi"
$splitChar = "s"

# Locate the LinkedIn URL in the document body.
$full = $d.Content.Text
$urlStart = $full.IndexOf($oldUrl)

if ($urlStart -ge 0) {
    $suffixStart = $urlStart + $keepPrefix.Length + $splitChar.Length
    $urlEnd = $urlStart + $oldUrl.Length

    # Remove the trailing "-799214175/" part of the old URL, leaving
    # "https://www.linkedin.com/in/jordan-alexis" in a single run.
    $suffixRange = $d.Range($suffixStart, $urlEnd)
    $suffixRange.Text = ""

    # Split the trailing "s" of "alexis" into its own run (matching the
    # target document, which stores it as a separate run with identical
    # formatting) by nudging a character property on just that range.
    $sRange = $d.Range($suffixStart - $splitChar.Length, $suffixStart)
    $sRange.Font.Bold = $true
    $sRange.Font.Bold = $false
}
